# change log output 2018-3-27-14-27:00
#
# Appends the latest exchange-rate snapshot (dated 2018-03-27 13:51:00) as a
# new row on the aggregate "allData_sheet" and as a new row on each
# individual currency sheet (CNY/JPY/GBP/EUR/RUB), without touching any of
# the previously logged rows.

$wb = $excel.ActiveWorkbook

# Helper: write a value as genuine text (not auto-coerced to a number),
# while leaving the cell's style back at the default "Normal" so no new
# cell formatting is left behind.
function Set-TextValue($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$updateDate = "2018-03-27 13:51:00"

$rates = @(
    @{ Name = "CNY"; Rate = "6.2604" },
    @{ Name = "JPY"; Rate = "105.6540" },
    @{ Name = "GBP"; Rate = "0.7029" },
    @{ Name = "EUR"; Rate = "0.8033" },
    @{ Name = "RUB"; Rate = "57.1192" }
)

# 1) Append one row per currency to the aggregate sheet.
$wsAll = $wb.Worksheets.Item("allData_sheet")
$nextRow = $wsAll.Cells.Item($wsAll.Rows.Count, 1).End(-4162).Row + 1

foreach ($entry in $rates) {
    $addrA = "A" + $nextRow
    $addrB = "B" + $nextRow
    $addrC = "C" + $nextRow
    Set-TextValue $wsAll $addrA $entry.Name
    Set-TextValue $wsAll $addrB $entry.Rate
    Set-TextValue $wsAll $addrC $updateDate
    $nextRow = $nextRow + 1
}

# 2) Append the same snapshot as a new row on each currency's own sheet.
foreach ($entry in $rates) {
    $ws = $wb.Worksheets.Item($entry.Name)
    $row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
    $addrA2 = "A" + $row
    $addrB2 = "B" + $row
    Set-TextValue $ws $addrA2 $entry.Rate
    Set-TextValue $ws $addrB2 $updateDate
}

Write-Output "Appended 2018-03-27 13:51:00 exchange-rate snapshot."
